$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the relevant paragraphs robustly via Find (text is unique),
# then resolve the containing paragraph index by scanning the
# Paragraphs collection for the one whose range covers the hit.
# ------------------------------------------------------------------
function Get-ParagraphIndexContaining([int]$start, [int]$end) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Start -le $start -and $cand.Range.End -ge $end) {
            return $i
        }
    }
    return -1
}

$wsWmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- find "Fecha de culminacion: ${fechaFin}" paragraph ---
$hit = $d.Content
$hit.Find.Execute("Fecha de culminaci", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fechaParaIdx = Get-ParagraphIndexContaining $hit.Start $hit.End

# the empty paragraph directly preceding it gets the new field content
$emptyBeforeIdx = $fechaParaIdx - 1

# --- find "Por la atencion" paragraph ---
$hit2 = $d.Content
$hit2.Find.Execute("Por la atenci", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$porParaIdx = Get-ParagraphIndexContaining $hit2.Start $hit2.End

$emptyAfterIdx = $porParaIdx + 1

# ------------------------------------------------------------------
# 1) Remove the empty paragraph that follows "Por la atencion..."
#    (delete the later one first so earlier indices stay valid)
# ------------------------------------------------------------------
$d.Paragraphs.Item($emptyAfterIdx).Range.Delete()

# ------------------------------------------------------------------
# 2) Remove the empty paragraph that precedes "Por la atencion..."
#    (its formatting is identical to the surviving paragraph's own
#    pPr, so the merge is a no-op formatting-wise)
# ------------------------------------------------------------------
$d.Paragraphs.Item($porParaIdx - 1).Range.Delete()

# ------------------------------------------------------------------
# 3) Clear the "Fecha de culminacion: ${fechaFin}" paragraph down to
#    an empty paragraph, and drop the Bold/BoldCs from its paragraph
#    mark formatting.
# ------------------------------------------------------------------
$fechaPara = $d.Paragraphs.Item($fechaParaIdx)
$emptyXml = '<w:p ' + $wsWmlNs + '><w:pPr><w:pStyle w:val="Standard"/><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-EC"/></w:rPr></w:pPr></w:p>'
$fechaPara.Range.InsertXML($emptyXml)

# ------------------------------------------------------------------
# 4) Populate the paragraph that used to be empty (right before the
#    one we just cleared) with the "Fecha de culminacion: " label and
#    a FILLIN field for ${fechas} (renamed from ${fechaFin}).
# ------------------------------------------------------------------
$newFieldXml = '<w:p ' + $wsWmlNs + '>' + `
    '<w:pPr><w:pStyle w:val="Standard"/><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-EC"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-EC"/></w:rPr><w:t xml:space="preserve">Fecha de culminaci&#243;n: </w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-EC"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-EC"/></w:rPr><w:instrText xml:space="preserve"> FILLIN  ${fechas}  \* MERGEFORMAT </w:instrText></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-EC"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-EC"/></w:rPr><w:t>${fechas}</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-EC"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' + `
    '</w:p>'

$d.Paragraphs.Item($emptyBeforeIdx).Range.InsertXML($newFieldXml)

Write-Host "done"
